$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$address = "https://github.com/ankurnecessary/dsa/blob/main/1_binarySearch/3_find_a_peak_element.js"
$title   = "dsa/3_find_a_peak_element.js at main · ankurnecessary/dsa · GitHub"

# Add the new hyperlink for F5. Passing $address as both the Address and the
# TextToDisplay makes Excel record the target URL in the <hyperlink display="..">
# attribute (matching the pattern already used by F3/F4); we then overwrite the
# cell's visible text with the friendly title, same as the existing rows.
$hl = $ws.Hyperlinks.Add($ws.Range("F5"), $address, "", "", $address)
$ws.Range("F5").Value = $title

# Match the existing hyperlink-cell formatting (style index used by F4) instead
# of letting Excel mint a brand-new style for F5.
$ws.Range("F4").Copy()
$ws.Range("F5").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Move the active selection to B6, as in the target workbook.
$ws.Range("B6").Select() | Out-Null
